# Insert a new data row at row 382 (pushing existing rows 382:506 down to 383:507)
# and populate it with the new record described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one whole row above the current row 382; Excel automatically shifts
# everything below (including formatting) down by one row.
$ws.Rows("382:382").Insert()

# Fill in the newly inserted row 382 with the new record's data.
$ws.Range("A382").Value = 1
$ws.Range("B382").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C382").Value = 'Arica y Parinacota'
$ws.Range("D382").Value = 44985
$ws.Range("E382").Value = 15
$ws.Range("F382").Value = 100112023
$ws.Range("G382").Value = 'Brócoli'
$ws.Range("H382").Value = 'Sin especificar'
$ws.Range("I382").Value = 'Tercera'
$ws.Range("J382").Value = 1200
$ws.Range("K382").Value = 400
$ws.Range("L382").Value = 500
$ws.Range("M382").Value = 450
$ws.Range("N382").Value = '$/unidad'
$ws.Range("O382").Value = 'Región de Arica y Parinacota'
$ws.Range("P382").Value = 450
$ws.Range("Q382").Value = 1
$ws.Range("R382").Value = 'Hortaliza'
